$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 1159.6154
$ws.Cells.Item(96, 9).Value = 1016.625
$ws.Cells.Item(96, 10).Value = 1388.4
$ws.Cells.Item(96, 11).Value = 3049.875
$ws.Cells.Item(96, 12).Value = 4165.200000000001
$ws.Cells.Item(96, 13).Value = -1676.875
$ws.Cells.Item(96, 14).Value = -6911.200000000001

$ws.Cells.Item(116, 8).Value = 114367.63
$ws.Cells.Item(116, 9).Value = 153106.08
$ws.Cells.Item(116, 10).Value = 5900
$ws.Cells.Item(116, 11).Value = 153106.08
$ws.Cells.Item(116, 12).Value = 5900
$ws.Cells.Item(116, 13).Value = -149664.08
$ws.Cells.Item(116, 14).Value = -12784

$ws.Cells.Item(121, 8).Value = 1149.7028
$ws.Cells.Item(121, 9).Value = 560
$ws.Cells.Item(121, 10).Value = 1201.7354
$ws.Cells.Item(121, 11).Value = 1680
$ws.Cells.Item(121, 12).Value = 3605.2062
$ws.Cells.Item(121, 13).Value = 67
$ws.Cells.Item(121, 14).Value = -7099.206200000001

$ws.Cells.Item(137, 8).Value = 3664.5356
$ws.Cells.Item(137, 9).Value = 6149.5
$ws.Cells.Item(137, 10).Value = 2284
$ws.Cells.Item(137, 11).Value = 18448.5
$ws.Cells.Item(137, 12).Value = 6852
$ws.Cells.Item(137, 13).Value = -15898.5
$ws.Cells.Item(137, 14).Value = -11952

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 734
$ws.Cells.Item(2, 9).Value = 660.2222
$ws.Cells.Item(2, 10).Value = 900
$ws.Cells.Item(2, 11).Value = 660.2222
$ws.Cells.Item(2, 12).Value = 900
$ws.Cells.Item(2, 13).Value = -547.2222
$ws.Cells.Item(2, 14).Value = -1126

$ws.Cells.Item(32, 8).Value = 8229.195
$ws.Cells.Item(32, 9).Value = 4733.579
$ws.Cells.Item(32, 10).Value = 52507
$ws.Cells.Item(32, 11).Value = 4733.579
$ws.Cells.Item(32, 12).Value = 52507
$ws.Cells.Item(32, 13).Value = -4446.579
$ws.Cells.Item(32, 14).Value = -53081

$ws.Cells.Item(45, 8).Value = 1248.6666
$ws.Cells.Item(45, 9).Value = 1143.5
$ws.Cells.Item(45, 10).Value = 1388.8889
$ws.Cells.Item(45, 11).Value = 1143.5
$ws.Cells.Item(45, 12).Value = 1388.8889
$ws.Cells.Item(45, 13).Value = -766.5
$ws.Cells.Item(45, 14).Value = -2142.8889

$ws.Cells.Item(74, 8).Value = 2713.5454
$ws.Cells.Item(74, 9).Value = 1600.091
$ws.Cells.Item(74, 11).Value = 1600.091
$ws.Cells.Item(74, 13).Value = -726.0909999999999

$ws.Cells.Item(77, 8).Value = 2713.5454
$ws.Cells.Item(77, 9).Value = 1600.091
$ws.Cells.Item(77, 11).Value = 8000.455
$ws.Cells.Item(77, 13).Value = -3632.455

$ws.Cells.Item(116, 8).Value = 734
$ws.Cells.Item(116, 9).Value = 660.2222
$ws.Cells.Item(116, 10).Value = 900
$ws.Cells.Item(116, 11).Value = 660.2222
$ws.Cells.Item(116, 12).Value = 900
$ws.Cells.Item(116, 13).Value = 1633.7778
$ws.Cells.Item(116, 14).Value = -5488

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 734
$ws.Cells.Item(3, 9).Value = 660.2222
$ws.Cells.Item(3, 10).Value = 900
$ws.Cells.Item(3, 11).Value = 660.2222
$ws.Cells.Item(3, 12).Value = 900
$ws.Cells.Item(3, 13).Value = -546.2222
$ws.Cells.Item(3, 14).Value = -1128

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1942.1613
$ws.Cells.Item(132, 9).Value = 1135.5897
$ws.Cells.Item(132, 11).Value = 3406.7691
$ws.Cells.Item(132, 13).Value = -876.7691

$ws.Cells.Item(134, 8).Value = 1522.9231
$ws.Cells.Item(134, 9).Value = 1084
$ws.Cells.Item(134, 10).Value = 2714.2856
$ws.Cells.Item(134, 11).Value = 3252
$ws.Cells.Item(134, 12).Value = 8142.8568
$ws.Cells.Item(134, 13).Value = -717
$ws.Cells.Item(134, 14).Value = -13212.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1688.8889
$ws.Cells.Item(75, 9).Value = 500
$ws.Cells.Item(75, 10).Value = 1837.5
$ws.Cells.Item(75, 11).Value = 1500
$ws.Cells.Item(75, 12).Value = 5512.5
$ws.Cells.Item(75, 13).Value = -502
$ws.Cells.Item(75, 14).Value = -7508.5

$ws.Cells.Item(78, 8).Value = 1688.8889
$ws.Cells.Item(78, 9).Value = 500
$ws.Cells.Item(78, 10).Value = 1837.5
$ws.Cells.Item(78, 11).Value = 4500
$ws.Cells.Item(78, 12).Value = 16537.5
$ws.Cells.Item(78, 13).Value = 492
$ws.Cells.Item(78, 14).Value = -26521.5

$ws.Cells.Item(107, 8).Value = 347.1591
$ws.Cells.Item(107, 9).Value = 224.06897
$ws.Cells.Item(107, 10).Value = 585.13336
$ws.Cells.Item(107, 11).Value = 672.20691
$ws.Cells.Item(107, 12).Value = 1755.40008
$ws.Cells.Item(107, 13).Value = 1247.79309
$ws.Cells.Item(107, 14).Value = -5595.40008

$ws.Cells.Item(129, 8).Value = 2481.8462
$ws.Cells.Item(129, 9).Value = 2493.5
$ws.Cells.Item(129, 10).Value = 2463.2
$ws.Cells.Item(129, 11).Value = 7480.5
$ws.Cells.Item(129, 12).Value = 7389.599999999999
$ws.Cells.Item(129, 13).Value = -2480.5
$ws.Cells.Item(129, 14).Value = -17389.6

$ws.Cells.Item(131, 8).Value = 2034.9215
$ws.Cells.Item(131, 9).Value = 2294.889
$ws.Cells.Item(131, 10).Value = 1893.1212
$ws.Cells.Item(131, 11).Value = 6884.667
$ws.Cells.Item(131, 12).Value = 5679.363600000001
$ws.Cells.Item(131, 13).Value = -1844.667
$ws.Cells.Item(131, 14).Value = -15759.3636

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 872.6842
$ws.Cells.Item(97, 9).Value = 411.66666
$ws.Cells.Item(97, 10).Value = 1663
$ws.Cells.Item(97, 11).Value = 411.66666
$ws.Cells.Item(97, 12).Value = 1663
$ws.Cells.Item(97, 13).Value = 84.33334000000002
$ws.Cells.Item(97, 14).Value = -2655

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1527
$ws.Cells.Item(61, 9).Value = 1406.9166
$ws.Cells.Item(61, 10).Value = 2247.5
$ws.Cells.Item(61, 11).Value = 1406.9166
$ws.Cells.Item(61, 12).Value = 2247.5
$ws.Cells.Item(61, 13).Value = -1204.9166
$ws.Cells.Item(61, 14).Value = -2651.5

$ws.Cells.Item(100, 8).Value = 32262290
$ws.Cells.Item(100, 9).Value = 5086.15
$ws.Cells.Item(100, 10).Value = 90911750
$ws.Cells.Item(100, 11).Value = 5086.15
$ws.Cells.Item(100, 12).Value = 90911750
$ws.Cells.Item(100, 13).Value = -4545.15
$ws.Cells.Item(100, 14).Value = -90912832

$ws.Cells.Item(113, 8).Value = 1527
$ws.Cells.Item(113, 9).Value = 1406.9166
$ws.Cells.Item(113, 10).Value = 2247.5
$ws.Cells.Item(113, 11).Value = 1406.9166
$ws.Cells.Item(113, 12).Value = 2247.5
$ws.Cells.Item(113, 13).Value = 763.0834
$ws.Cells.Item(113, 14).Value = -6587.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 12949.5
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 12949.5
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 12949.5
$ws.Cells.Item(96, 13).ClearContents()
$ws.Cells.Item(96, 14).Value = -15695.5

$ws.Cells.Item(100, 8).Value = 1529.5883
$ws.Cells.Item(100, 9).Value = 1562.5
$ws.Cells.Item(100, 10).Value = 1003
$ws.Cells.Item(100, 11).Value = 3125
$ws.Cells.Item(100, 12).Value = 2006
$ws.Cells.Item(100, 13).Value = -2584
$ws.Cells.Item(100, 14).Value = -3088
